$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, preventing Excel from
# auto-converting numeric-looking strings (e.g. "212.93") into numbers.
function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.242.96"
$ws.Range("E2").Value = "  +3.66%  "
$ws.Range("D3").Value = "1.607.25"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("E4").Value = "  -0.13%  "
Set-TextCell "D5" "212.93"
$ws.Range("E5").Value = "  +2.56%  "
Set-TextCell "D7" "0.486"
$ws.Range("E7").Value = "  +2.06%  "
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("D12").Value = "1.833.50"
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("D13").Value = "1.608.12"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("E14").Value = "  -0.38%  "
Set-TextCell "D15" "0.513"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "26.236.02"
$ws.Range("E16").Value = "  +3.69%  "
Set-TextCell "D17" "60.84"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("E18").Value = "  +2.16%  "
Set-TextCell "D19" "208.75"
$ws.Range("E19").Value = "  +12.45%  "
$ws.Range("E20").Value = "  -0.21%  "
Set-TextCell "D21" "4.26"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("E22").Value = "  +0.47%  "
Set-TextCell "D23" "6.02"
$ws.Range("E23").Value = "  +2.26%  "
Set-TextCell "D24" "1.82"
$ws.Range("E24").Value = "  +10.02%  "
Set-TextCell "D25" "142.13"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("E26").Value = "  -0.16%  "
Set-TextCell "D27" "0.125"
$ws.Range("E27").Value = "  -4.06%  "
Set-TextCell "D28" "15.27"
$ws.Range("E28").Value = "  +2.87%  "
Set-TextCell "D29" "6.46"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("E32").Value = "  +3.11%  "
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("D36").Value = "1.110.15"
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("E37").Value = "  +6.96%  "
$ws.Range("E38").Value = "  +0.17%  "
Set-TextCell "D39" "2.33"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("E40").Value = "  +1.31%  "
Set-TextCell "D41" "0.497"
$ws.Range("E41").Value = "  +0.19%  "
Set-TextCell "D42" "0.776"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "1.745.95"
$ws.Range("E43").Value = "  +2.71%  "
Set-TextCell "D44" "92.90"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "0.0₆0107"
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("E47").Value = "  +9.13%  "
Set-TextCell "D48" "53.61"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E50").Value = "  +0.91%  "
